# Adds 4 new expense rows (24 julio) to the "Gastos" sheet, rows 10-13,
# mirroring the author's manual data entry in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gastos")

# Helper: write a value that *looks* numeric/currency (e.g. "$500") while
# keeping it as literal text, matching how the source file stores amounts
# like "$413.52" as plain strings rather than formatted numbers.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 10: Desarmador de estrella - $500 (Herramientas)
$ws.Range("A10").Value = "24  julio"
$ws.Range("B10").Value = "Equipo y Mob"
$ws.Range("C10").Value = "Desarmador de estrella"
Set-TextValue $ws.Range("D10") '$500'
$ws.Range("E10").Value = "Herramientas"

# Row 11: Martillo Nuevo - $150 (Herramientas)
$ws.Range("A11").Value = "24  julio"
$ws.Range("B11").Value = "Equipo y Mob"
$ws.Range("C11").Value = "Martillo Nuevo"
Set-TextValue $ws.Range("D11") '$150'
$ws.Range("E11").Value = "Herramientas"

# Row 12: Equipo y Mobiliario - $  (Herramientas)
$ws.Range("A12").Value = "24  julio"
$ws.Range("B12").Value = "Equipo y Mob"
$ws.Range("C12").Value = "Equipo y Mobiliario"
$ws.Range("D12").Value = "$ "
$ws.Range("E12").Value = "Herramientas"

# Row 13: Compra en Transito - $ (Publicidad)
$ws.Range("A13").Value = "24  julio"
$ws.Range("B13").Value = "Compra en Transito"
$ws.Range("C13").Value = "Compra en Transito"
$ws.Range("D13").Value = "$"
$ws.Range("E13").Value = "Publicidad"
